# Progress update as of 04-Nov-2025
#
# 1) "Training Dashboard": every training's "PERIOD TO EXPIRE" (col H) drops
#    by one day and "LAST UPDATE" (col I) moves from 03-Nov-2025 to
#    04-Nov-2025, for all data rows (3-18).
# 2) "Exam Dashboard": a new exam result ("Cs Hoist") is inserted as row 6
#    (pushing the old "TOTAL AVERAGE" summary row down to row 7) and the
#    average is recalculated to 74.28%.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Training Dashboard
# ---------------------------------------------------------------------
$training = $wb.Worksheets.Item("Training Dashboard")

for ($r = 3; $r -le 18; $r++) {
    $period = $training.Cells.Item($r, 8).Value2
    $training.Cells.Item($r, 8).Value = $period - 1
    # Leading "'" keeps this a literal text value (matches the source
    # file, where LAST UPDATE is stored as text, not a real date).
    $training.Cells.Item($r, 9).Value = "'04-Nov-2025"
}

# ---------------------------------------------------------------------
# 2) Exam Dashboard
# ---------------------------------------------------------------------
$exam = $wb.Worksheets.Item("Exam Dashboard")

# Insert a new row above the old "TOTAL AVERAGE" row (row 6), shifting it
# to row 7, and copy the formatting of the previous data row (5) onto it.
$exam.Rows(6).Insert()
$exam.Range("A5:G5").Copy()
$exam.Range("A6:G6").PasteSpecial(-4122)

$exam.Range("A6").Value = 4
$exam.Range("B6").Value = "Cs Hoist"
$exam.Range("C6").Value = "'30-Oct-2025"
$exam.Range("D6").Value = "'71.61%"
$exam.Range("E6").Value = "low percentage"
$exam.Range("F6").Value = "This is a low mark, please retake the exam and improve your score. date is valid"

# Recalculated TOTAL AVERAGE on the (now shifted) summary row.
$exam.Range("D7").Value = "'74.28%"
